$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.181.89'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.658.81'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = "'219.24"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = "'0.5231"
$ws.Range("E6").Value = '  -1.50%  '
$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("D8").Value = "'0.2629"
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = "'0.06296"
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").Value = "'20.58"
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").Value = "'0.07818"
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = "'4.492"
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '1.656.70'
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '1.887.52'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = "'0.5543"
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '0.0₅8024'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").Value = "'65.14"
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '26.205.33'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = "'4.633"
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").Value = "'196.28"
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("E22").Value = '  -0.64%  '
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = "'0.1200"
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").Value = "'7.141"
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = "'1.495"
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("D30").Value = "'0.05743"
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("D31").Value = "'1.276"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = "'3.491"
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("D33").Value = "'3.367"
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").Value = "'1.584"
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'0.9548"
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = "'2.807"
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = "'0.5720"
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = "'5.976"
$ws.Range("E40").Value = '  +2.35%  '
$ws.Range("D41").Value = '1.060.49'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").Value = "'0.8474"
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = "'103.91"
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '1.798.34'
$ws.Range("D46").Value = "'58.19"
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").Value = "'1.011"
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₈105'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.4406"
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").Value = "'8.048"
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").Value = "'0.05202"
$ws.Range("E51").Value = '  +0.77%  '

# Remove the quote-prefix style artifact so these cells keep the default style
# (matches the original workbook, which stored them as plain inline strings).
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
